# Adding pickling to replication
# Fill in the standard-error rows (row 4 = theta_se, row 6 = lambda_se)
# which previously held placeholder "(nan)" text, with the actual
# computed (bootstrapped) standard errors, one value per hour column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row4 = @{
    "B4" = "(0.54)"
    "C4" = "(0.05)"
    "D4" = "(0.09)"
    "E4" = "(1.3)"
    "F4" = "(0.01)"
    "G4" = "(1.79)"
    "H4" = "(1.11)"
    "I4" = "(1.92)"
    "J4" = "(4.24)"
}

$row6 = @{
    "B6" = "(0.81)"
    "C6" = "(0.19)"
    "D6" = "(0.42)"
    "E6" = "(0.43)"
    "F6" = "(0.2)"
    "G6" = "(1.41)"
    "H6" = "(0.32)"
    "I6" = "(0.77)"
    "J6" = "(2.88)"
}

$columns = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

foreach ($col in $columns) {
    $ws.Range("$col" + "4").Value = $row4["$col" + "4"]
    $ws.Range("$col" + "6").Value = $row6["$col" + "6"]
}
